# Process 2 (Main) Updated
# - Small Bug Fixes
# - Added Extensive Reporting
# - Added More Log Messages
# - Added Extra Try Catch, to reduce errors

$wb = $excel.ActiveWorkbook

# --- Update the Customer_Database sheet: shift every date in column I
#     (rows 2-51) forward by one day (e.g. 6/29/2020 -> 6/30/2020).
$ws = $wb.Worksheets.Item("Customer_Database")

for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 9)   # column I
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}

# --- Update the window/selection state on the Customer_Database sheet
#     (this also drops the old "topLeftCell" scroll anchor, matching the
#     saved view after the user scrolled/selected a new cell).
$ws.Activate()
$ws.Range("M14").Select()

# --- Update the workbook window size/position (best effort - mirrors the
#     saved <workbookView> bounds after the user resized/moved the window).
try {
    $win = $wb.Windows.Item(1)
    $win.Left = -18705
    $win.Top = 7815
    $win.Width = 17700
    $win.Height = 10380
} catch {
    # window geometry not settable in this host; ignore
}

$excel.Left = -18705
$excel.Top = 7815
$excel.Width = 17700
$excel.Height = 10380
